# ------------------------------------------------------------------
# scripts/data/info.xlsx -- "update process and scripts"
#
# Adds a "fullname" column (D) giving the full binomial/latin name of
# each organism, a "mechanism of action" column (B) for each
# antibiotic, tidies up the "neither Gram+/-" wording, freezes the
# first column, and pads a few trailing blank rows that Excel keeps
# formatted for future entries.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("B1").Font.Name = "Arial"
$ws.Range("B1").Font.Size = 10
$ws.Range("D1").Value = "fullname"

# --- D column: latin binomial ("fullname") for each bacterium -------
# (E. coli filled first to mirror how the source sheet was populated)
$ws.Range("D13").Value = "Escherichia coli"
$ws.Range("D13").Font.Name = "Arial"
$ws.Range("D13").Font.Size = 10
$ws.Range("D2").Value = "Acetinobacter baumannii"
$ws.Range("D2").Font.Name = "Arial"
$ws.Range("D2").Font.Size = 10
$ws.Range("D3").Value = "Klebsiella pneumoniae"
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 10
$ws.Range("D4").Value = "Enterococcus faecium"
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 10
$ws.Range("D5").Value = "Neisseria gonorrhoeae"
$ws.Range("D5").Font.Name = "Arial"
$ws.Range("D5").Font.Size = 10
$ws.Range("D7").Value = "Mycobacterium tuberculosis"
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("D7").Font.Size = 10
$ws.Range("D8").Value = "Proteus mirabilis"
$ws.Range("D8").Font.Name = "Arial"
$ws.Range("D8").Font.Size = 10
$ws.Range("D10").Value = "Clostridium difficile"
$ws.Range("D10").Font.Name = "Arial"
$ws.Range("D10").Font.Size = 10
$ws.Range("D11").Value = "Pseudomonas aeruginosa"
$ws.Range("D11").Font.Name = "Arial"
$ws.Range("D11").Font.Size = 10
$ws.Range("D12").Value = "Streptococcus pneumoniae"
$ws.Range("D12").Font.Name = "Arial"
$ws.Range("D12").Font.Size = 10
$ws.Range("D14").Value = "Staphylococcus aureus"
$ws.Range("D14").Font.Name = "Arial"
$ws.Range("D14").Font.Size = 10

# D6 (Shigella) and D9 (CoNS) reuse text already present elsewhere
$ws.Range("D6").Value = "Shigella"
$ws.Range("D6").Font.Name = "Arial"
$ws.Range("D6").Font.Size = 10
$ws.Range("D9").Value = "Coagulase-negative Staphilococci"
$ws.Range("D9").Font.Name = "Arial"
$ws.Range("D9").Font.Size = 10

# --- C7: reword the "neither Gram+/-" category -----------------------
$ws.Range("C7").Value = "Neither Gram-negative nor Gram-positive"

# --- B column: antibiotic mechanism of action, rows 15-35 ------------
$ws.Range("B15").Value = "Cell wall"
$ws.Range("B15").Font.Name = "Arial"
$ws.Range("B15").Font.Size = 10
$ws.Range("B16").Value = "Protein synthesis"
$ws.Range("B16").Font.Name = "Arial"
$ws.Range("B16").Font.Size = 10
$ws.Range("B27").Value = "DNA synthesis"
$ws.Range("B27").Font.Name = "Arial"
$ws.Range("B27").Font.Size = 10
$ws.Range("B24").Value = "Metabolism"
$ws.Range("B24").Font.Name = "Arial"
$ws.Range("B24").Font.Size = 10
$ws.Range("B17").Value = "Protein synthesis"
$ws.Range("B17").Font.Name = "Arial"
$ws.Range("B17").Font.Size = 10
$ws.Range("B18").Value = "Cell wall"
$ws.Range("B18").Font.Name = "Arial"
$ws.Range("B18").Font.Size = 10
$ws.Range("B19").Value = "Cell wall"
$ws.Range("B19").Font.Name = "Arial"
$ws.Range("B19").Font.Size = 10
$ws.Range("B20").Value = "Protein synthesis"
$ws.Range("B20").Font.Name = "Arial"
$ws.Range("B20").Font.Size = 10
$ws.Range("B21").Value = "Cell wall"
$ws.Range("B21").Font.Name = "Arial"
$ws.Range("B21").Font.Size = 10
$ws.Range("B22").Value = "Cell wall"
$ws.Range("B22").Font.Name = "Arial"
$ws.Range("B22").Font.Size = 10
$ws.Range("B23").Value = "Cell wall"
$ws.Range("B23").Font.Name = "Arial"
$ws.Range("B23").Font.Size = 10
$ws.Range("B25").Value = "Cell wall"
$ws.Range("B25").Font.Name = "Arial"
$ws.Range("B25").Font.Size = 10
$ws.Range("B26").Value = "Protein synthesis"
$ws.Range("B26").Font.Name = "Arial"
$ws.Range("B26").Font.Size = 10
$ws.Range("B28").Value = "Protein synthesis"
$ws.Range("B28").Font.Name = "Arial"
$ws.Range("B28").Font.Size = 10
$ws.Range("B29").Value = "Metabolism"
$ws.Range("B29").Font.Name = "Arial"
$ws.Range("B29").Font.Size = 10
$ws.Range("B30").Value = "Cell wall"
$ws.Range("B30").Font.Name = "Arial"
$ws.Range("B30").Font.Size = 10
$ws.Range("B31").Value = "Cell wall"
$ws.Range("B31").Font.Name = "Arial"
$ws.Range("B31").Font.Size = 10
$ws.Range("B32").Value = "Cell wall"
$ws.Range("B32").Font.Name = "Arial"
$ws.Range("B32").Font.Size = 10
$ws.Range("B33").Value = "DNA synthesis"
$ws.Range("B33").Font.Name = "Arial"
$ws.Range("B33").Font.Size = 10
$ws.Range("B34").Value = "DNA synthesis"
$ws.Range("B34").Font.Name = "Arial"
$ws.Range("B34").Font.Size = 10
$ws.Range("B35").Value = "Cell wall"
$ws.Range("B35").Font.Name = "Arial"
$ws.Range("B35").Font.Size = 10

# --- D column placeholder (styled, value-less) for antibiotic rows ---
$ws.Range("D15").Font.Name = "Arial"
$ws.Range("D15").Font.Size = 10
$ws.Range("D16").Font.Name = "Arial"
$ws.Range("D16").Font.Size = 10
$ws.Range("D17").Font.Name = "Arial"
$ws.Range("D17").Font.Size = 10
$ws.Range("D18").Font.Name = "Arial"
$ws.Range("D18").Font.Size = 10
$ws.Range("D19").Font.Name = "Arial"
$ws.Range("D19").Font.Size = 10
$ws.Range("D20").Font.Name = "Arial"
$ws.Range("D20").Font.Size = 10
$ws.Range("D21").Font.Name = "Arial"
$ws.Range("D21").Font.Size = 10
$ws.Range("D22").Font.Name = "Arial"
$ws.Range("D22").Font.Size = 10
$ws.Range("D23").Font.Name = "Arial"
$ws.Range("D23").Font.Size = 10
$ws.Range("D24").Font.Name = "Arial"
$ws.Range("D24").Font.Size = 10
$ws.Range("D25").Font.Name = "Arial"
$ws.Range("D25").Font.Size = 10
$ws.Range("D26").Font.Name = "Arial"
$ws.Range("D26").Font.Size = 10
$ws.Range("D27").Font.Name = "Arial"
$ws.Range("D27").Font.Size = 10
$ws.Range("D28").Font.Name = "Arial"
$ws.Range("D28").Font.Size = 10
$ws.Range("D29").Font.Name = "Arial"
$ws.Range("D29").Font.Size = 10
$ws.Range("D30").Font.Name = "Arial"
$ws.Range("D30").Font.Size = 10
$ws.Range("D31").Font.Name = "Arial"
$ws.Range("D31").Font.Size = 10
$ws.Range("D32").Font.Name = "Arial"
$ws.Range("D32").Font.Size = 10
$ws.Range("D33").Font.Name = "Arial"
$ws.Range("D33").Font.Size = 10
$ws.Range("D34").Font.Name = "Arial"
$ws.Range("D34").Font.Size = 10
$ws.Range("D35").Font.Name = "Arial"
$ws.Range("D35").Font.Size = 10

# --- Trailing blank rows 36-41 (kept styled for future entries) ------
$ws.Range("A36").Font.Name = "Arial"
$ws.Range("A36").Font.Size = 10
$ws.Range("D36").Font.Name = "Arial"
$ws.Range("D36").Font.Size = 10
$ws.Range("A37").Font.Name = "Arial"
$ws.Range("A37").Font.Size = 10
$ws.Range("D37").Font.Name = "Arial"
$ws.Range("D37").Font.Size = 10
$ws.Range("A38").Font.Name = "Arial"
$ws.Range("A38").Font.Size = 10
$ws.Range("D38").Font.Name = "Arial"
$ws.Range("D38").Font.Size = 10
$ws.Range("A39").Font.Name = "Arial"
$ws.Range("A39").Font.Size = 10
$ws.Range("D39").Font.Name = "Arial"
$ws.Range("D39").Font.Size = 10
$ws.Range("A40").Font.Name = "Arial"
$ws.Range("A40").Font.Size = 10
$ws.Range("D40").Font.Name = "Arial"
$ws.Range("D40").Font.Size = 10
$ws.Range("A41").Font.Name = "Arial"
$ws.Range("A41").Font.Size = 10
$ws.Range("D41").Font.Name = "Arial"
$ws.Range("D41").Font.Size = 10

# --- View: freeze first column, leave selection on D22 in right pane -
$ws.Activate()
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D22").Select()
